$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "AQ2" = 41
    "G2" = 1.62
    "I2" = 5.25
    "Q2" = 1.62
    "R2" = 2.3
    "S2" = 2.03
    "T2" = 1.83
    "Y2" = 1.67
    "Z2" = 2.1
    "Y3" = 1.7
    "AG4" = 4.65
    "AI4" = 18.5
    "AK4" = 101
    "AL4" = 5.4
    "AM4" = 10
    "AN4" = 10
    "AO4" = 27
    "AP4" = 27
    "AQ4" = 50
    "H4" = 2.6
    "K4" = 1.78
    "L4" = 3.15
    "M4" = 1.16
    "N4" = 4.65
    "O4" = 1.62
    "P4" = 2.15
    "Q4" = 2.82
    "R4" = 1.38
    "U4" = 5.1
    "V4" = 1.13
    "W4" = 1.65
    "X4" = 2.12
    "Y4" = 2.18
    "Z4" = 1.6
    "AA5" = 8.25
    "AB5" = 19.5
    "AC5" = 13.5
    "AE5" = 45
    "AF5" = 60
    "AG5" = 4.9
    "AH5" = 5.4
    "AI5" = 17.5
    "AK5" = 101
    "AL5" = 5.3
    "AN5" = 9.5
    "AO5" = 23
    "AP5" = 24
    "AQ5" = 45
    "G5" = 3.9
    "H5" = 2.67
    "I5" = 2.22
    "J5" = 4.45
    "K5" = 1.83
    "L5" = 2.95
    "N5" = 4.9
    "O5" = 1.57
    "P5" = 2.27
    "Q5" = 2.65
    "R5" = 1.42
    "U5" = 4.8
    "W5" = 1.6
    "X5" = 2.22
    "AA6" = 5.8
    "AB6" = 7.5
    "AD6" = 13.5
    "AE6" = 15
    "AG6" = 6.2
    "AH6" = 6.4
    "AI6" = 17.5
    "AJ6" = 100
    "AL6" = 11
    "AM6" = 28
    "AN6" = 17
    "AO6" = 100
    "G6" = 1.72
    "H6" = 3.25
    "I6" = 5.2
    "J6" = 2.25
    "K6" = 2.07
    "L6" = 5.4
    "M6" = 1.09
    "N6" = 6.2
    "O6" = 1.39
    "P6" = 2.77
    "Q6" = 2.15
    "R6" = 1.62
    "U6" = 3.65
    "V6" = 1.24
    "W6" = 1.44
    "X6" = 2.6
    "Y6" = 2
    "Z6" = 1.72
    "J7" = 2.38
    "M7" = 1.07
    "O7" = 1.36
    "V7" = 1.25
    "J8" = 1.8
    "K8" = 2.38
    "M8" = 1.05
    "O8" = 1.3
    "Q8" = 1.91
    "R8" = 1.91
    "V8" = 1.33
    "N10" = 9
    "Q10" = 2.1
    "R10" = 1.73
    "AO11" = 29
    "J11" = 3.25
    "W13" = 1.33
    "W14" = 1.22
    "Q15" = 1.75
    "W15" = 1.3
    "AA16" = 9
    "AB16" = 11
    "AD16" = 21
    "AE16" = 17
    "AG16" = 12
    "AN16" = 11
    "AP16" = 23
    "G16" = 2.15
    "H16" = 3.5
    "I16" = 3
    "J16" = 2.88
    "K16" = 2.2
    "L16" = 3.6
    "M16" = 1.05
    "N16" = 11
    "Q16" = 1.88
    "R16" = 1.98
    "U16" = 3.25
    "V16" = 1.33
    "W16" = 1.33
    "Q17" = 1.62
    "W17" = 1.27
    "AA19" = 6.9
    "AD19" = 10.75
    "AE19" = 12
    "AF19" = 25
    "AG19" = 11
    "AH19" = 7.6
    "AI19" = 17
    "AL19" = 15.5
    "AM19" = 37
    "AN19" = 18.5
    "AO19" = 120
    "AP19" = 65
    "AQ19" = 60
    "G19" = 1.52
    "H19" = 3.85
    "I19" = 5.8
    "J19" = 2.05
    "K19" = 2.18
    "L19" = 5.7
    "O19" = 1.24
    "P19" = 3.3
    "Q19" = 1.72
    "R19" = 1.88
    "U19" = 2.67
    "V19" = 1.36
    "AA26" = 4.4
    "AB26" = 4.8
    "AC26" = 6.8
    "AD26" = 7.6
    "AF26" = 30
    "AH26" = 6.2
    "AI26" = 19
    "AL26" = 11
    "AM26" = 30
    "AN26" = 18
    "AP26" = 80
    "AQ26" = 80
    "G26" = 1.48
    "H26" = 3.95
    "I26" = 6.6
    "J26" = 1.98
    "K26" = 2.24
    "L26" = 6.8
    "V26" = 1.27
    "Y26" = 2.22
    "Z26" = 1.61
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
